# Applies the "release build, some tweaks" edits described by the diff.
#
# Strategy: this runtime's Range.Text / Find.Execute assignment normalizes
# (merges) every run inside the touched paragraph into a single run, which
# does not match the target OOXML (which keeps/creates specific run and
# <w:proofErr/> boundaries). So instead we replace each target paragraph's
# whole range with a hand-built <w:p>...</w:p> fragment via Range.InsertXML,
# which is applied verbatim (run/proofErr boundaries preserved exactly).

$d = $word.ActiveDocument
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

function Set-ParagraphXml($paraIndex, $innerRunsXml) {
    $p = $d.Paragraphs($paraIndex)
    $r = $p.Range
    $xml = "<w:p $wNs>" + $innerRunsXml + "</w:p>"
    $r.InsertXML($xml)
}

# --- Hunk 1 (was paragraph 22): merge "really" + " have to multiply..." runs
$apos = [char]0x2019
$inner22 = '<w:r><w:t xml:space="preserve">Since we split the number up in multiples of 10' + $apos + 's, you </w:t></w:r>'
$inner22 += '<w:r><w:t>only</w:t></w:r>'
$inner22 += '<w:r><w:t xml:space="preserve"> </w:t></w:r>'
$inner22 += '<w:r><w:t>really have to multiply the non-zero numbers.</w:t></w:r>'
Set-ParagraphXml 22 $inner22

# --- Hunk 2 (was paragraph 53): split "mult" out with proofErr markers
$inner53 = '<w:r><w:t xml:space="preserve">[double digit </w:t></w:r>'
$inner53 += '<w:proofErr w:type="spellStart"/>'
$inner53 += '<w:r><w:t>mult</w:t></w:r>'
$inner53 += '<w:proofErr w:type="spellEnd"/>'
$inner53 += '<w:r><w:t>. Area illustration]</w:t></w:r>'
Set-ParagraphXml 53 $inner53

# --- Hunk 3 (was paragraph 59): split "mult" out with proofErr markers
$inner59 = '<w:r><w:t>[evaluate phase</w:t></w:r>'
$inner59 += '<w:r><w:t xml:space="preserve">, wait for double digit </w:t></w:r>'
$inner59 += '<w:proofErr w:type="spellStart"/>'
$inner59 += '<w:r><w:t>mult</w:t></w:r>'
$inner59 += '<w:proofErr w:type="spellEnd"/>'
$inner59 += '<w:r><w:t>]</w:t></w:r>'
Set-ParagraphXml 59 $inner59

# --- Hunk 4 (was paragraph 67): "our dimension" -> "the sky"
$inner67 = '<w:r><w:t xml:space="preserve">Multiple space blobs have pierced through </w:t></w:r>'
$inner67 += '<w:r><w:t>the sky</w:t></w:r>'
$inner67 += '<w:r><w:t>!</w:t></w:r>'
Set-ParagraphXml 67 $inner67

# --- Hunk 5a (was paragraph 69): split into 3 runs, append ", and wreak havoc"
$inner69 = '<w:r><w:t>We must banish them immediately before they fall down to Earth</w:t></w:r>'
$inner69 += '<w:r><w:t>, and wreak havoc</w:t></w:r>'
$inner69 += '<w:r><w:t>!</w:t></w:r>'
Set-ParagraphXml 69 $inner69

# --- Hunk 5b (was paragraph 70): split "blobology" out with proofErr markers
$inner70 = '<w:r><w:t xml:space="preserve">With our latest advancements in </w:t></w:r>'
$inner70 += '<w:proofErr w:type="spellStart"/>'
$inner70 += '<w:r><w:t>blobology</w:t></w:r>'
$inner70 += '<w:proofErr w:type="spellEnd"/>'
$inner70 += '<w:r><w:t>, we will be deploying Attack Blobs.</w:t></w:r>'
Set-ParagraphXml 70 $inner70

# --- Hunk 6: append a blank paragraph, then two new paragraphs at the very
# end of the document body (before sectPr).
$end = $d.Content.End
$r = $d.Range($end, $end)
$r.InsertXML("<w:p $wNs/>")

$n = $d.Paragraphs.Count
$p = $d.Paragraphs($n)
$p.Range.InsertParagraphAfter()

$n = $d.Paragraphs.Count
$p = $d.Paragraphs($n)
$xmlNew1 = "<w:p $wNs><w:r><w:t>Multiple space blobs have pierced through the sky! Banish these threats with the power of mathematics before they wreak havoc!</w:t></w:r></w:p>"
$p.Range.InsertXML($xmlNew1)

$n = $d.Paragraphs.Count
$p = $d.Paragraphs($n)
$p.Range.InsertParagraphAfter()

$n = $d.Paragraphs.Count
$p = $d.Paragraphs($n)
$inner76 = '<w:r><w:t>Use the touchpad or mouse to drag a</w:t></w:r>'
$inner76 += '<w:r><w:t>nd connect a</w:t></w:r>'
$inner76 += '<w:r><w:t xml:space="preserve"> blob to another. Once connected, you will go through the process of generating the product of the equation. If correct, an attack blob will </w:t></w:r>'
$inner76 += '<w:r><w:t>appear to banish the paired blobs</w:t></w:r>'
$inner76 += '<w:r><w:t>.</w:t></w:r>'
$xmlNew2 = "<w:p $wNs>" + $inner76 + "</w:p>"
$p.Range.InsertXML($xmlNew2)

Write-Output "done"
